$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: econ_program_prop_smearacf
$ws.Range("A20").Value = "econ_program_prop_smearacf"
$ws.Range("B20").Value = "yes"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "yes"
$ws.Range("E20").Value = 0
$ws.Range("BE20").Value = 0
$ws.Range("BM20").Value = 0

# Row 21: econ_program_totalcost_smearacf
$ws.Range("A21").Value = "econ_program_totalcost_smearacf"
$ws.Range("B21").Value = "yes"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "yes"
$ws.Range("E21").Value = 0
$ws.Range("AO21").Value = 0
$ws.Range("BE21").Value = 0
$ws.Range("BM21").Value = 0

# Row 22: econ_program_reflectioncost_smearacf
$ws.Range("A22").Value = "econ_program_reflectioncost_smearacf"
$ws.Range("B22").Value = "yes"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "yes"
$ws.Range("E22").Value = 0
$ws.Range("AH22").Value = 0
$ws.Range("BD22").Value = 0
$ws.Range("BM22").Value = 0

# Row 23: econ_program_unitcost_smearacf
$ws.Range("A23").Value = "econ_program_unitcost_smearacf"
$ws.Range("B23").Value = "yes"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "yes"
$ws.Range("E23").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("AL23").Value = 0
$ws.Range("BD23").Value = 0
$ws.Range("BM23").Value = 0

# Row 24: econ_program_prop_xpertacf
$ws.Range("A24").Value = "econ_program_prop_xpertacf"
$ws.Range("B24").Value = "yes"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = "yes"
$ws.Range("E24").Value = 0
$ws.Range("BB24").Value = 0
$ws.Range("BM24").Value = 0

# Row 25: econ_program_totalcost_xpertacf
$ws.Range("A25").Value = "econ_program_totalcost_xpertacf"
$ws.Range("B25").Value = "yes"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "yes"
$ws.Range("E25").Value = 0
$ws.Range("BB25").Value = 0
$ws.Range("BM25").Value = 0

# Row 26: econ_program_reflectioncost_xpertacf
$ws.Range("A26").Value = "econ_program_reflectioncost_xpertacf"
$ws.Range("B26").Value = "yes"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "yes"
$ws.Range("E26").Value = 0
$ws.Range("BB26").Value = 0
$ws.Range("BM26").Value = 0

# Row 27: econ_program_unitcost_xpertacf
$ws.Range("A27").Value = "econ_program_unitcost_xpertacf"
$ws.Range("B27").Value = "yes"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "yes"
$ws.Range("E27").Value = 0
$ws.Range("BB27").Value = 0
$ws.Range("BM27").Value = 0

# Page setup: paper size 9 (A4), portrait orientation.
# (side effect: legacyDrawing relationship id is renumbered, matching target)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("BB6").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BB23").Select()
